# Commit: "Updated qotd style, added cards to carousel for portfolio,
#          added bootstrap icon"
#
# The OOXML diff shows <w:strike/> being added to four ListParagraph
# bullets (to both the paragraph mark's rPr and every run's rPr, i.e. the
# whole paragraph gets struck through) -- the same pattern already used
# elsewhere in the doc (e.g. "Properly use a dynamic root path...").
#
# Target bullets (matched by their visible text, not just position, so the
# script fails loudly instead of silently striking the wrong line if the
# document is reshuffled):
#   - "Use Bootstrap icons"
#   - "Use Bootstrap and style to position elements for a professional
#      looking website"
#   - "Ignore backup files and folders"
#   - "Remove ignored files and folders from the repository"

$d = $word.ActiveDocument

$targets = @(
    "Use Bootstrap icons",
    "Use Bootstrap and style to position elements for a professional looking website",
    "Ignore backup files and folders",
    "Remove ignored files and folders from the repository"
)

foreach ($needle in $targets) {
    $found = $false
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $t = $p.Range.Text.Trim()
        if ($t -eq $needle) {
            # Strike the whole paragraph - the paragraph-mark rPr as well
            # as every run in it - matching the diff exactly.
            $p.Range.Font.StrikeThrough = $true
            $found = $true
            break
        }
    }
    if (-not $found) {
        throw "Could not locate paragraph with text: $needle"
    }
}
